$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 93 <- data from former row 99 (id 6236254)
$ws.Range("B93").Value = 6236254
$ws.Range("C93").Value = 'Venezuela Primera Division'
$ws.Range("D93").Value = 'Venezuela Primera Division'
$ws.Range("E93").Value = 45199.6875
$ws.Range("F93").Value = 'Academia Puerto Cabello'
$ws.Range("G93").Value = 'Estudiantes Merida'
$ws.Range("H93").Value = 1
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 'H'
$ws.Range("K93").Value = 1.727
$ws.Range("L93").Value = 3.4
$ws.Range("M93").Value = 4.333
$ws.Range("N93").Value = 1.666
$ws.Range("O93").Value = 3.4
$ws.Range("P93").Value = 4.75
$ws.Range("Q93").Value = -0.75
$ws.Range("R93").Value = 1.875
$ws.Range("S93").Value = 1.925
$ws.Range("T93").Value = 2.5
$ws.Range("U93").Value = 1.9
$ws.Range("V93").Value = 1.9
$ws.Range("W93").Value = 0.6659999999999999
$ws.Range("X93").Value = -1
$ws.Range("Y93").Value = -1
$ws.Range("Z93").Value = 0.4375
$ws.Range("AA93").Value = -0.5
$ws.Range("AB93").Value = -1
$ws.Range("AC93").Value = 0.8999999999999999

# Row 94 <- data from former row 93 (id 6236251)
$ws.Range("B94").Value = 6236251
$ws.Range("C94").Value = 'Venezuela Primera Division'
$ws.Range("D94").Value = 'Venezuela Primera Division'
$ws.Range("E94").Value = 45199.6875
$ws.Range("F94").Value = 'Angostura FC'
$ws.Range("G94").Value = 'Portuguesa'
$ws.Range("H94").Value = 1
$ws.Range("I94").Value = 2
$ws.Range("J94").Value = 'A'
$ws.Range("K94").Value = 3.1
$ws.Range("L94").Value = 3.2
$ws.Range("M94").Value = 2.15
$ws.Range("N94").Value = 4
$ws.Range("O94").Value = 3.6
$ws.Range("P94").Value = 1.75
$ws.Range("Q94").Value = 0.75
$ws.Range("R94").Value = 1.8
$ws.Range("S94").Value = 2
$ws.Range("T94").Value = 2.5
$ws.Range("U94").Value = 1.95
$ws.Range("V94").Value = 1.85
$ws.Range("W94").Value = -1
$ws.Range("X94").Value = -1
$ws.Range("Y94").Value = 0.75
$ws.Range("Z94").Value = -0.5
$ws.Range("AA94").Value = 0.5
$ws.Range("AB94").Value = 0.95
$ws.Range("AC94").Value = -1

# Row 96 <- data from former row 98 (id 6236253)
$ws.Range("B96").Value = 6236253
$ws.Range("C96").Value = 'Venezuela Primera Division'
$ws.Range("D96").Value = 'Venezuela Primera Division'
$ws.Range("E96").Value = 45199.6875
$ws.Range("F96").Value = 'Deportivo La Guaira'
$ws.Range("G96").Value = 'UCV'
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 'D'
$ws.Range("K96").Value = 1.833
$ws.Range("L96").Value = 3.25
$ws.Range("M96").Value = 4
$ws.Range("N96").Value = 2
$ws.Range("O96").Value = 3.2
$ws.Range("P96").Value = 3.5
$ws.Range("Q96").Value = -0.25
$ws.Range("R96").Value = 1.775
$ws.Range("S96").Value = 2.025
$ws.Range("T96").Value = 2.25
$ws.Range("U96").Value = 1.9
$ws.Range("V96").Value = 1.9
$ws.Range("W96").Value = -1
$ws.Range("X96").Value = 2.2
$ws.Range("Y96").Value = -1
$ws.Range("Z96").Value = -0.5
$ws.Range("AA96").Value = 0.5125
$ws.Range("AB96").Value = -1
$ws.Range("AC96").Value = 0.8999999999999999

# Row 97 <- data from former row 94 (id 6236611)
$ws.Range("B97").Value = 6236611
$ws.Range("C97").Value = 'Venezuela Primera Division'
$ws.Range("D97").Value = 'Venezuela Primera Division'
$ws.Range("E97").Value = 45199.6875
$ws.Range("F97").Value = 'Mineros'
$ws.Range("G97").Value = 'Monagas'
$ws.Range("H97").Value = 2
$ws.Range("I97").Value = 1
$ws.Range("J97").Value = 'H'
$ws.Range("K97").Value = 3.2
$ws.Range("L97").Value = 3.4
$ws.Range("M97").Value = 2
$ws.Range("N97").Value = 4.2
$ws.Range("O97").Value = 3.8
$ws.Range("P97").Value = 1.65
$ws.Range("Q97").Value = 0.75
$ws.Range("R97").Value = 1.95
$ws.Range("S97").Value = 1.85
$ws.Range("T97").Value = 2.5
$ws.Range("U97").Value = 1.825
$ws.Range("V97").Value = 1.975
$ws.Range("W97").Value = 3.2
$ws.Range("X97").Value = -1
$ws.Range("Y97").Value = -1
$ws.Range("Z97").Value = 0.95
$ws.Range("AA97").Value = -1
$ws.Range("AB97").Value = 0.825
$ws.Range("AC97").Value = -1

# Row 98 <- data from former row 97 (id 6236252)
$ws.Range("B98").Value = 6236252
$ws.Range("C98").Value = 'Venezuela Primera Division'
$ws.Range("D98").Value = 'Venezuela Primera Division'
$ws.Range("E98").Value = 45199.6875
$ws.Range("F98").Value = 'Deportivo Tachira'
$ws.Range("G98").Value = 'CD Hermanos Colmenares'
$ws.Range("H98").Value = 1
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 'H'
$ws.Range("K98").Value = 1.363
$ws.Range("L98").Value = 4.2
$ws.Range("M98").Value = 7.5
$ws.Range("N98").Value = 1.333
$ws.Range("O98").Value = 4.5
$ws.Range("P98").Value = 8
$ws.Range("Q98").Value = -1.5
$ws.Range("R98").Value = 2
$ws.Range("S98").Value = 1.8
$ws.Range("T98").Value = 2.5
$ws.Range("U98").Value = 1.925
$ws.Range("V98").Value = 1.875
$ws.Range("W98").Value = 0.333
$ws.Range("X98").Value = -1
$ws.Range("Y98").Value = -1
$ws.Range("Z98").Value = -1
$ws.Range("AA98").Value = 0.8
$ws.Range("AB98").Value = -1
$ws.Range("AC98").Value = 0.875

# Row 99 <- data from former row 96 (id 6236255)
$ws.Range("B99").Value = 6236255
$ws.Range("C99").Value = 'Venezuela Primera Division'
$ws.Range("D99").Value = 'Venezuela Primera Division'
$ws.Range("E99").Value = 45199.6875
$ws.Range("F99").Value = 'Deportivo Rayo Zuliano'
$ws.Range("G99").Value = 'Caracas'
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 'D'
$ws.Range("K99").Value = 3.75
$ws.Range("L99").Value = 3.1
$ws.Range("M99").Value = 1.95
$ws.Range("N99").Value = 2.9
$ws.Range("O99").Value = 2.875
$ws.Range("P99").Value = 2.45
$ws.Range("Q99").Value = 0.25
$ws.Range("R99").Value = 1.775
$ws.Range("S99").Value = 2.025
$ws.Range("T99").Value = 2.25
$ws.Range("U99").Value = 1.85
$ws.Range("V99").Value = 1.95
$ws.Range("W99").Value = -1
$ws.Range("X99").Value = 1.875
$ws.Range("Y99").Value = -1
$ws.Range("Z99").Value = 0.3875
$ws.Range("AA99").Value = -0.5
$ws.Range("AB99").Value = -1
$ws.Range("AC99").Value = 0.95

# Row 102 <- data from former row 103 (id 6236615)
$ws.Range("B102").Value = 6236615
$ws.Range("C102").Value = 'Venezuela Primera Division'
$ws.Range("D102").Value = 'Venezuela Primera Division'
$ws.Range("E102").Value = 45206.6875
$ws.Range("F102").Value = 'Deportivo Rayo Zuliano'
$ws.Range("G102").Value = 'Academia Puerto Cabello'
$ws.Range("H102").Value = 1
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 'H'
$ws.Range("K102").Value = 2.375
$ws.Range("L102").Value = 3.3
$ws.Range("M102").Value = 2.625
$ws.Range("N102").Value = 2.45
$ws.Range("O102").Value = 3.2
$ws.Range("P102").Value = 2.55
$ws.Range("Q102").Value = 0
$ws.Range("R102").Value = 1.875
$ws.Range("S102").Value = 1.925
$ws.Range("T102").Value = 2.5
$ws.Range("U102").Value = 2
$ws.Range("V102").Value = 1.8
$ws.Range("W102").Value = 1.45
$ws.Range("X102").Value = -1
$ws.Range("Y102").Value = -1
$ws.Range("Z102").Value = 0.875
$ws.Range("AA102").Value = -1
$ws.Range("AB102").Value = -1
$ws.Range("AC102").Value = 0.8

# Row 103 <- data from former row 102 (id 6236616)
$ws.Range("B103").Value = 6236616
$ws.Range("C103").Value = 'Venezuela Primera Division'
$ws.Range("D103").Value = 'Venezuela Primera Division'
$ws.Range("E103").Value = 45206.6875
$ws.Range("F103").Value = 'UCV'
$ws.Range("G103").Value = 'Metropolitanos FC'
$ws.Range("H103").Value = 3
$ws.Range("I103").Value = 2
$ws.Range("J103").Value = 'H'
$ws.Range("K103").Value = 3.3
$ws.Range("L103").Value = 3.2
$ws.Range("M103").Value = 2.05
$ws.Range("N103").Value = 2.75
$ws.Range("O103").Value = 3.2
$ws.Range("P103").Value = 2.3
$ws.Range("Q103").Value = 0.25
$ws.Range("R103").Value = 1.75
$ws.Range("S103").Value = 2.05
$ws.Range("T103").Value = 2.5
$ws.Range("U103").Value = 1.975
$ws.Range("V103").Value = 1.825
$ws.Range("W103").Value = 1.75
$ws.Range("X103").Value = -1
$ws.Range("Y103").Value = -1
$ws.Range("Z103").Value = 0.75
$ws.Range("AA103").Value = -1
$ws.Range("AB103").Value = 0.9750000000000001
$ws.Range("AC103").Value = -1

# Row 135 <- data from former row 136 (id 7842507)
$ws.Range("B135").Value = 7842507
$ws.Range("C135").Value = 'Venezuela Primera Division'
$ws.Range("D135").Value = 'Venezuela Primera Division'
$ws.Range("E135").Value = 45339.78125
$ws.Range("F135").Value = 'Academia Puerto Cabello'
$ws.Range("G135").Value = 'Estudiantes Merida'
$ws.Range("H135").Value = 2
$ws.Range("I135").Value = 1
$ws.Range("J135").Value = 'H'
$ws.Range("K135").Value = 1.727
$ws.Range("L135").Value = 3.5
$ws.Range("M135").Value = 4.2
$ws.Range("N135").Value = 1.85
$ws.Range("O135").Value = 3.5
$ws.Range("P135").Value = 3.6
$ws.Range("Q135").Value = -0.5
$ws.Range("R135").Value = 1.925
$ws.Range("S135").Value = 1.875
$ws.Range("T135").Value = 2.5
$ws.Range("U135").Value = 1.9
$ws.Range("V135").Value = 1.9
$ws.Range("W135").Value = 0.8500000000000001
$ws.Range("X135").Value = -1
$ws.Range("Y135").Value = -1
$ws.Range("Z135").Value = 0.925
$ws.Range("AA135").Value = -1
$ws.Range("AB135").Value = 0.8999999999999999
$ws.Range("AC135").Value = -1

# Row 136 <- data from former row 135 (id 7842504)
$ws.Range("B136").Value = 7842504
$ws.Range("C136").Value = 'Venezuela Primera Division'
$ws.Range("D136").Value = 'Venezuela Primera Division'
$ws.Range("E136").Value = 45339.78125
$ws.Range("F136").Value = 'Angostura FC'
$ws.Range("G136").Value = 'Deportivo La Guaira'
$ws.Range("H136").Value = 1
$ws.Range("I136").Value = 1
$ws.Range("J136").Value = 'D'
$ws.Range("K136").Value = 2.75
$ws.Range("L136").Value = 3
$ws.Range("M136").Value = 2.45
$ws.Range("N136").Value = 3.1
$ws.Range("O136").Value = 2.875
$ws.Range("P136").Value = 2.3
$ws.Range("Q136").Value = 0.25
$ws.Range("R136").Value = 1.8
$ws.Range("S136").Value = 2
$ws.Range("T136").Value = 2.25
$ws.Range("U136").Value = 2.05
$ws.Range("V136").Value = 1.75
$ws.Range("W136").Value = -1
$ws.Range("X136").Value = 1.875
$ws.Range("Y136").Value = -1
$ws.Range("Z136").Value = 0.4
$ws.Range("AA136").Value = -0.5
$ws.Range("AB136").Value = -0.5
$ws.Range("AC136").Value = 0.375

# Row 162 <- data from former row 163 (id 7952893)
$ws.Range("B162").Value = 7952893
$ws.Range("C162").Value = 'Venezuela Primera Division'
$ws.Range("D162").Value = 'Venezuela Primera Division'
$ws.Range("E162").Value = 45366.83333333334
$ws.Range("F162").Value = 'UCV'
$ws.Range("G162").Value = 'Deportivo La Guaira'
$ws.Range("H162").Value = 1
$ws.Range("I162").Value = 1
$ws.Range("J162").Value = 'D'
$ws.Range("K162").Value = 2.1
$ws.Range("L162").Value = 3
$ws.Range("M162").Value = 3.25
$ws.Range("N162").Value = 2.25
$ws.Range("O162").Value = 3.1
$ws.Range("P162").Value = 2.9
$ws.Range("Q162").Value = -0.25
$ws.Range("R162").Value = 2.025
$ws.Range("S162").Value = 1.775
$ws.Range("T162").Value = 2
$ws.Range("U162").Value = 1.8
$ws.Range("V162").Value = 2
$ws.Range("W162").Value = -1
$ws.Range("X162").Value = 2.1
$ws.Range("Y162").Value = -1
$ws.Range("Z162").Value = -0.5
$ws.Range("AA162").Value = 0.3875
$ws.Range("AB162").Value = 0
# AC162 unchanged (stays -0)

# Row 163 <- data from former row 162 (id 7952905)
$ws.Range("B163").Value = 7952905
$ws.Range("C163").Value = 'Venezuela Primera Division'
$ws.Range("D163").Value = 'Venezuela Primera Division'
$ws.Range("E163").Value = 45366.83333333334
$ws.Range("F163").Value = 'Angostura FC'
$ws.Range("G163").Value = 'Deportivo Tachira'
$ws.Range("H163").Value = 2
$ws.Range("I163").Value = 0
$ws.Range("J163").Value = 'H'
$ws.Range("K163").Value = 3.6
$ws.Range("L163").Value = 3.6
$ws.Range("M163").Value = 1.8
$ws.Range("N163").Value = 3.75
$ws.Range("O163").Value = 2.875
$ws.Range("P163").Value = 2.1
$ws.Range("Q163").Value = 0.25
$ws.Range("R163").Value = 1.95
$ws.Range("S163").Value = 1.85
$ws.Range("T163").Value = 2
$ws.Range("U163").Value = 2.025
$ws.Range("V163").Value = 1.775
$ws.Range("W163").Value = 2.75
$ws.Range("X163").Value = -1
$ws.Range("Y163").Value = -1
$ws.Range("Z163").Value = 0.95
$ws.Range("AA163").Value = -1
$ws.Range("AB163").Value = 0
# AC163 unchanged (stays -0)
